$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H76" = 5793935.5
    "I76" = 171999.67
    "J76" = 7199419.5
    "K76" = 171999.67
    "L76" = 7199419.5
    "M76" = -171684.67
    "N76" = -7200049.5
    "H79" = 5793935.5
    "I79" = 171999.67
    "J79" = 7199419.5
    "K79" = 171999.67
    "L79" = 7199419.5
    "M79" = -170907.67
    "N79" = -7201603.5
    "H92" = 95005.14
    "J92" = 221240.89
    "L92" = 221240.89
    "N92" = -223736.89
    "H116" = 38959104
    "I116" = 27892106
    "J116" = 55559610
    "K116" = 27892106
    "L116" = 55559610
    "M116" = -27888664
    "N116" = -55566494
    "H127" = 1116.2084
    "I127" = 1110.85
    "J127" = 1143
    "K127" = 3332.55
    "L127" = 3429
    "M127" = 1627.45
    "N127" = -13349
    "H131" = 8230.532999999999
    "I131" = 2619.75
    "K131" = 7859.25
    "M131" = -2819.25
    "H137" = 26317802
    "J137" = 2464.95
    "L137" = 7394.849999999999
    "N137" = -12494.85
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H2" = 690.6
    "J2" = 0
    "L2" = 0
    "H32" = 3791.0625
    "I32" = 4508.6
    "J32" = 1228.4286
    "K32" = 4508.6
    "L32" = 1228.4286
    "M32" = -4221.6
    "N32" = -1802.4286
    "H45" = 1378.4286
    "I45" = 1124.25
    "J45" = 1717.3334
    "K45" = 1124.25
    "L45" = 1717.3334
    "M45" = -747.25
    "N45" = -2471.3334
    "H46" = 4224.6
    "I46" = 2936.75
    "J46" = 5083.1665
    "K46" = 2936.75
    "L46" = 5083.1665
    "M46" = -2617.75
    "N46" = -5721.1665
    "H61" = 1481.2307
    "I61" = 1282.8
    "K61" = 1282.8
    "M61" = -1070.8
    "H74" = 5231.8423
    "I74" = 5588.25
    "K74" = 5588.25
    "M74" = -4714.25
    "H77" = 5231.8423
    "I77" = 5588.25
    "K77" = 27941.25
    "M77" = -23573.25
    "H92" = 49999
    "J92" = 49999
    "L92" = 49999
    "N92" = -54991
    "H97" = 604.7742
    "I97" = 619.96155
    "J97" = 525.8
    "K97" = 619.96155
    "L97" = 525.8
    "M97" = -123.96155
    "N97" = -1517.8
    "H116" = 690.6
    "J116" = 0
    "L116" = 0
    "H132" = 1314.9
    "I132" = 1213.5714
    "J132" = 1551.3334
    "K132" = 3640.7142
    "L132" = 4654.0002
    "M132" = -1110.7142
    "N132" = -9714.0002
    "H136" = 1481.2307
    "I136" = 1282.8
    "K136" = 3848.4
    "M136" = -1298.4
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
$clears = @("N2", "N116")
foreach ($key in $clears) {
    $ws.Range($key).ClearContents()
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H3" = 690.6
    "J3" = 0
    "L3" = 0
    "H9" = 20000
    "J9" = 20000
    "L9" = 20000
    "N9" = -20336
    "H37" = 868.625
    "I37" = 868.625
    "J37" = 0
    "K37" = 868.625
    "L37" = 0
    "N37" = -731.625
    "H94" = 3134.7856
    "I94" = 3308.6667
    "J94" = 2613.1428
    "K94" = 3308.6667
    "L94" = 2613.1428
    "M94" = -2857.6667
    "N94" = -3515.1428
    "H135" = 49900
    "J135" = 49900
    "L135" = 49900
    "N135" = -60040
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
$clears = @("N3", "M37")
foreach ($key in $clears) {
    $ws.Range($key).ClearContents()
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H7" = 84.13333
    "I7" = 84.13333
    "J7" = 0
    "K7" = 84.13333
    "L7" = 0
    "N7" = 28.86667
    "H16" = 2740.5715
    "J16" = 4128.3335
    "L16" = 4128.3335
    "N16" = -4702.3335
    "H31" = 5917.1875
    "J31" = 3364.8108
    "L31" = 3364.8108
    "N31" = -3954.8108
    "H34" = 5917.1875
    "J34" = 3364.8108
    "L34" = 3364.8108
    "N34" = -3768.8108
    "H58" = 1415.8605
    "I58" = 1340
    "K58" = 1340
    "M58" = -1137
    "H113" = 2740.5715
    "J113" = 4128.3335
    "L113" = 4128.3335
    "N113" = -8468.333500000001
    "H122" = 2082.2368
    "I122" = 2148.9644
    "J122" = 1895.4
    "K122" = 6446.8932
    "L122" = 5686.200000000001
    "M122" = -3996.8932
    "N122" = -10586.2
    "H134" = 1894
    "I134" = 1766.0435
    "K134" = 5298.1305
    "M134" = -2763.1305
    "H136" = 1415.8605
    "I136" = 1340
    "K136" = 4020
    "M136" = -1470
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
$clears = @("M7")
foreach ($key in $clears) {
    $ws.Range($key).ClearContents()
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H122" = 739.75
    "J122" = 685
    "L122" = 6165
    "N122" = -11065
    "H131" = 5062
    "I131" = 773.5
    "J131" = 8360.846
    "K131" = 2320.5
    "L131" = 25082.538
    "M131" = 2719.5
    "N131" = -35162.538
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H122" = 2148.3845
    "I122" = 1831
    "K122" = 5493
    "M122" = -3043
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H82" = 2542.611
    "I82" = 2660.2354
    "J82" = 543
    "K82" = 2660.2354
    "L82" = 543
    "M82" = -2299.2354
    "N82" = -1265
    "H85" = 2542.611
    "I85" = 2660.2354
    "J85" = 543
    "K85" = 2660.2354
    "L85" = 543
    "M85" = -1412.2354
    "N85" = -3039
    "H132" = 4057
    "I132" = 2485.8667
    "K132" = 7457.6001
    "M132" = -4927.6001
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H100" = 7689.4287
    "I100" = 7689.4287
    "K100" = 15378.8574
    "M100" = -14837.8574
    "H122" = 1342.1892
    "I122" = 1164.7241
    "J122" = 1985.5
    "K122" = 3494.1723
    "L122" = 5956.5
    "M122" = -1044.1723
    "N122" = -10856.5
    "H132" = 3078.8367
    "I132" = 3009.7568
    "J132" = 3291.8333
    "K132" = 9029.270400000001
    "L132" = 9875.499899999999
    "M132" = -6499.270400000001
    "N132" = -14935.4999
}
foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
